# Apply "Add data for 2022-08-19" update:
#  - Rename sheet from "Through 2022-08-10" to "Through 2022-08-11"
#  - Update header label in I1 from "2022 (through 08-10)" to "2022 (through 08-11)"
#  - Update August 2022 value (I9) from 52 to 61
#  - Update Total 2022 value (I14) from 1022 to 1031

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename worksheet tab
$ws.Name = "Through 2022-08-11"

# Update the running-total column header text
$ws.Range("I1").Value = "2022 (through 08-11)"

# Update August figure for 2022
$ws.Range("I9").Value = 61

# Update Total figure for 2022
$ws.Range("I14").Value = 1031
